$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 21:52"

# 2. Swap the country names for rows 77/78 (Uzbekistan now sorts above Republica de Macedonia)
$ws.Range("A77").Value = "Uzbekistan"
$ws.Range("A78").Value = "Republica de Macedonia"

# 3. Update Estados Unidos stats (row 4)
$ws.Range("B4").Value = 525559
$ws.Range("C4").Value = 22683
$ws.Range("E4").Value = 476500
$ws.Range("G4").Value = 1557
$ws.Range("H4").Value = 20304

# 4. Update row 77 (Uzbekistan) stats
$ws.Range("B77").Value = 767
$ws.Range("C77").Value = 143
$ws.Range("D77").Value = 42
$ws.Range("E77").Value = 721
$ws.Range("F77").Value = 8
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 4

# 5. Update row 78 (Republica de Macedonia) stats
$ws.Range("B78").Value = 760
$ws.Range("C78").Value = 49
$ws.Range("D78").Value = 41
$ws.Range("E78").Value = 685
$ws.Range("F78").Value = 15
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 34
